$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Julio de 2020 a las 21:27"

# Update country rows: name (if changed) + statistics columns B:H

# Row 4: Estados Unidos
$ws.Range("B4").Value = 4211641
$ws.Range("C4").Value = 41323
$ws.Range("D4").Value = 1996023
$ws.Range("E4").Value = 2067661
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 608
$ws.Range("H4").Value = 147957

# Row 6: India
$ws.Range("B6").Value = 1337021
$ws.Range("C6").Value = 48891
$ws.Range("D6").Value = 850048
$ws.Range("E6").Value = 455568
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 760
$ws.Range("H6").Value = 31405

# Row 8: Sudafrica
$ws.Range("B8").Value = 421996
$ws.Range("C8").Value = 13944
$ws.Range("D8").Value = 245771
$ws.Range("E8").Value = 169882
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 250
$ws.Range("H8").Value = 6343

# Row 22: Francia
$ws.Range("B22").Value = 180528
$ws.Range("C22").Value = 1130
$ws.Range("D22").Value = 80815
$ws.Range("E22").Value = 69521
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = 30192

# Row 24: Canada
$ws.Range("B24").Value = 113040
$ws.Range("C24").Value = 368
$ws.Range("D24").Value = 98766
$ws.Range("E24").Value = 5396
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 4
$ws.Range("H24").Value = 8878

# Row 65: Uzbekistan
$ws.Range("B65").Value = 19360
$ws.Range("C65").Value = 492
$ws.Range("D65").Value = 10472
$ws.Range("E65").Value = 8782
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 4
$ws.Range("H65").Value = 106

# Row 70: Costa de Marfil
$ws.Range("B70").Value = 15253
$ws.Range("C70").Value = 252
$ws.Range("D70").Value = 9499
$ws.Range("E70").Value = 5660
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 94

# Row 74: Costa Rica
$ws.Range("A74").Value = "Costa Rica"
$ws.Range("B74").Value = 13669
$ws.Range("C74").Value = 540
$ws.Range("D74").Value = 3505
$ws.Range("E74").Value = 10077
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 7
$ws.Range("H74").Value = 87

# Row 75: Venezuela
$ws.Range("A75").Value = "Venezuela"
$ws.Range("B75").Value = 13613
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 7752
$ws.Range("E75").Value = 5732
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 129

# Row 76: Australia
$ws.Range("A76").Value = "Australia"
$ws.Range("B76").Value = 13595
$ws.Range("C76").Value = 289
$ws.Range("D76").Value = 8775
$ws.Range("E76").Value = 4681
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 6
$ws.Range("H76").Value = 139

# Row 77: Dinamarca
$ws.Range("A77").Value = "Dinamarca"
$ws.Range("B77").Value = 13438
$ws.Range("C77").Value = 48
$ws.Range("D77").Value = 12340
$ws.Range("E77").Value = 485
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = 613

# Row 95: Mauritania
$ws.Range("B95").Value = 6116
$ws.Range("C95").Value = 49
$ws.Range("D95").Value = 4206
$ws.Range("E95").Value = 1754
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 156

# Row 118: Libia
$ws.Range("A118").Value = "Libia"
$ws.Range("B118").Value = 2424
$ws.Range("C118").Value = 110
$ws.Range("D118").Value = 504
$ws.Range("E118").Value = 1863
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 57

# Row 119: Hong Kong
$ws.Range("A119").Value = "Hong Kong"
$ws.Range("B119").Value = 2373
$ws.Range("C119").Value = 123
$ws.Range("D119").Value = 1407
$ws.Range("E119").Value = 950
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = 16

# Row 120: Zimbabue
$ws.Range("A120").Value = "Zimbabue"
$ws.Range("B120").Value = 2296
$ws.Range("C120").Value = 172
$ws.Range("D120").Value = 514
$ws.Range("E120").Value = 1750
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 4
$ws.Range("H120").Value = 32

# Row 121: Sudan del Sur
$ws.Range("A121").Value = "Sudan del Sur"
$ws.Range("B121").Value = 2258
$ws.Range("C121").Value = 19
$ws.Range("D121").Value = 1175
$ws.Range("E121").Value = 1038
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 45

# Row 122: Cabo Verde
$ws.Range("A122").Value = "Cabo Verde"
$ws.Range("B122").Value = 2220
$ws.Range("C122").Value = 30
$ws.Range("D122").Value = 1216
$ws.Range("E122").Value = 982
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = 22

# Row 150: Angola
$ws.Range("B150").Value = 880
$ws.Range("C150").Value = 29
$ws.Range("D150").Value = 241
$ws.Range("E150").Value = 604
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 2
$ws.Range("H150").Value = 35

# Row 151: Santo Tome y Principe
$ws.Range("A151").Value = "Santo Tome y Principe"
$ws.Range("B151").Value = 860
$ws.Range("C151").Value = 111
$ws.Range("D151").Value = 610
$ws.Range("E151").Value = 236
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 14

# Row 152: Togo
$ws.Range("A152").Value = "Togo"
$ws.Range("B152").Value = 828
$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 584
$ws.Range("E152").Value = 228
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 16

# Row 153: Jamaica
$ws.Range("A153").Value = "Jamaica"
$ws.Range("B153").Value = 821
$ws.Range("C153").Value = 5
$ws.Range("D153").Value = 711
$ws.Range("E153").Value = 100
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 10

# Row 210: Islas Malvinas
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("B210").Value = 13
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 13
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

# Row 211: Groenlandia
$ws.Range("A211").Value = "Groenlandia"
$ws.Range("B211").Value = 13
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 13
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0
